# Working data plots on single cascade
$wb = $excel.ActiveWorkbook

# --- Sheet "State Variables" (2nd sheet) ---
$wsState = $wb.Worksheets.Item(2)

# Row 2: replace single constant in C2 with repeated values across E2:W2
$wsState.Range("C2").ClearContents()
$wsState.Range("E2:W2").Value = 700

# Row 5: replace single constant in C5 with repeated values across E5:W5
$wsState.Range("C5").ClearContents()
$wsState.Range("E5:W5").Value = 1000

# Row 8: remove the constant in C8 (no replacement)
$wsState.Range("C8").ClearContents()

# --- Sheet "Parameters" (3rd sheet) keeps its selection at I22 but is no longer the active tab ---
$wsParams = $wb.Worksheets.Item(3)
$wsParams.Activate()
$wsParams.Range("I22").Select()

# Activate "State Variables" sheet last so it becomes the active/selected tab
$wsState.Activate()
$wsState.Range("M16").Select()
